$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Values are stored as text (percentages and numeric-looking strings),
# so force Text number format before assigning to avoid Excel auto-converting
# them into numeric/percentage values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.03%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.38%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.025"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.09%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07890"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.24%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.138"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-17.77%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.815"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.47%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.797"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.21%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9274"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.26%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1745"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.87%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07982"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.68%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08776"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.02%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03122"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.05%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1004"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.32%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001512"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.54%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005824"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.77%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.465"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.25%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.279"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.94%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3294"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.98%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.51%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.141"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.27%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1788"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.44%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04581"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.47%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001233"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.71%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004522"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.20%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.19%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01749"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.09%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04755"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.91%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007373"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.62%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1370"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.07%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002333"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "5.90%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01086"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "11.83%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006044"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.99%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.03%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.003390"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-61.21%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.8205"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.01%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002094"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.03%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001994"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.03%"
